$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9; this shifts existing rows 9..114 down to 10..115,
# preserving their values/styles exactly as required.
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with the new data values.
$ws.Range("D9").Value = 44503
$ws.Range("E9").Value = 4
$ws.Range("F9").Value = 100112021
$ws.Range("G9").Value = "Ají"
$ws.Range("H9").Value = "Inferno"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 560
$ws.Range("K9").Value = 29000
$ws.Range("L9").Value = 30000
$ws.Range("M9").Value = 29500
$ws.Range("N9").Value = "$/caja 12 kilos"
$ws.Range("O9").Value = "Región de Arica y Parinacota"
$ws.Range("P9").Value = 2458
$ws.Range("Q9").Value = 12
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "Terminal La Palmera de La Serena"
$ws.Range("C9").Value = "Coquimbo"
$ws.Range("R9").Value = "Hortaliza"
